$wb = $excel.ActiveWorkbook

# --- Insert new sheet "2022-Q1" between "2021-Q2" and "总计" ---
$ws2021 = $wb.Worksheets.Item("2021-Q2")
$newSheet = $wb.Worksheets.Add($null, $ws2021)
$newSheet.Name = "2022-Q1"

# Re-fetch sheet references AFTER all add/rename operations complete - references
# captured before a sheet insertion can end up pointing at the wrong sheet.
$ws2022 = $wb.Worksheets.Item("2022-Q1")
$wsTotal = $wb.Worksheets.Item("总计")

# ---------- Sheet "2022-Q1" ----------

# Header row values
$ws2022.Range("B1").Value = "基金代码"
$ws2022.Range("C1").Value = "基金名称"
$ws2022.Range("D1").Value = "基金规模"
$ws2022.Range("E1").Value = "股票总仓位"
$ws2022.Range("F1").Value = "仓位占比"
$ws2022.Range("G1").Value = "持有市值(亿元)"
$ws2022.Range("H1").Value = "仓位排名"
# Match header formatting to the rest of the workbook's summary-header style
$wsTotal.Range("B1").Copy()
$ws2022.Range("B1:H1").PasteSpecial(-4122)

# Data rows - column A (index numbers) and column H (rank numbers) are numeric;
# columns B-G hold numeric-looking identifiers/figures that must stay text.
$ws2022.Range("A2").Value = 0
$ws2022.Range("B2").Value = "'163110"
$ws2022.Range("C2").Value = "'申万菱信量化小盘股票(LOF)"
$ws2022.Range("D2").Value = "'5.68"
$ws2022.Range("E2").Value = "'92.25"
$ws2022.Range("F2").Value = "'1.35"
$ws2022.Range("G2").Value = "'0.0767"
$ws2022.Range("H2").Value = 3

$ws2022.Range("A3").Value = 1
$ws2022.Range("B3").Value = "'515860"
$ws2022.Range("C3").Value = "'嘉实中证新兴科技100策略ETF"
$ws2022.Range("D3").Value = "'2.25"
$ws2022.Range("E3").Value = "'98.94"
$ws2022.Range("F3").Value = "'2.17"
$ws2022.Range("G3").Value = "'0.0488"
$ws2022.Range("H3").Value = 10

$ws2022.Range("A4").Value = 2
$ws2022.Range("B4").Value = "'162907"
$ws2022.Range("C4").Value = "'泰信中证锐联基本面400指数（LOF）"
$ws2022.Range("D4").Value = "'0.23"
$ws2022.Range("E4").Value = "'94.61"
$ws2022.Range("F4").Value = "'0.54"
$ws2022.Range("G4").Value = "'0.0012"
$ws2022.Range("H4").Value = 10

# Column A shares the bold/bordered "index" style used throughout the workbook
$wsTotal.Range("A2").Copy()
$ws2022.Range("A2:A4").PasteSpecial(-4122)
$ws2022.Range("A2").Value = 0
$ws2022.Range("A3").Value = 1
$ws2022.Range("A4").Value = 2

# Strip the quote-prefix formatting the leading "'" added to B2:G4 by re-applying
# the plain (unstyled) format already used for text cells elsewhere in the workbook.
# (Values are already set above; pasting formats only - do not touch Value again,
# or the quote-prefix style would be re-applied.)
$ws2021.Range("B2").Copy()
$ws2022.Range("B2:G4").PasteSpecial(-4122)

# ---------- Sheet "总计" ----------
# Insert a new leading data row for the "2022-Q1" summary line, ahead of the
# existing "2021-Q2" row.
$wsTotal.Range("A2:D2").Insert()

# The inserted row inherits formatting from the row above (the bold header);
# re-apply the plain formatting already used by the existing data row instead.
$wsTotal.Range("A3:D3").Copy()
$wsTotal.Range("A2:D2").PasteSpecial(-4122)

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q1"
$wsTotal.Range("C2").Value = 3
$wsTotal.Range("D2").Value = 0.13
$wsTotal.Range("A3").Value = 1

# Restore the originally active sheet (adding "2022-Q1" shifted focus to it).
$ws2021.Activate()
